# Applies the StructureDefinition-insight-summary.xlsx update:
#  - Metadata sheet: Version 5.0.0 -> 6.0.0, Date bump, Publisher gets a value,
#    the two duplicate "Contact / No display for ContactDetail" rows collapse
#    into a single "Jurisdiction / United States of America" row, and the
#    Description text loses its trailing space.
#  - Elements sheet: the root Extension row's "Short"/"Definition" columns
#    (K2/L2) are updated from the generic "Extension" / "An Extension" to the
#    profile-specific "Insight Summary" / full description text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely,
# shifting Description and everything below it up by one row.
$meta.Rows.Item(11).Delete()

# Description text (now at row 11) loses its trailing space.
$meta.Range("B11").Value = "Summarizes the insight associated with the element that the insight summary extension is embedded in. The insight array element in the meta section of the resource with the same insight record id contains the full details on the insight."

# ---------------------------------------------------------------------------
# Elements sheet
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition columns change from the generic
# placeholder text to the profile-specific summary.
$elements.Range("K2").Value = "Insight Summary"
$elements.Range("L2").Value = "Summarizes the insight associated with the element that the insight summary extension is embedded in. The insight array element in the meta section of the resource with the same insight record id contains the full details on the insight."
